$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.935.53'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '1.643.24'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.53'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.66'
$ws.Range('E8').Value = '  +1.77%  '
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').Value = '1.875.79'
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('D13').Value = '1.653.16'
$ws.Range('E13').Value = '  +1.73%  '
$ws.Range('E14').Value = '  +4.23%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.77'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '27.908.18'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '230.45'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('D19').Value = '0.0₃0726'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.94'
$ws.Range('E22').Value = '  +4.91%  '
$ws.Range('E23').Value = '  +1.55%  '
$ws.Range('E24').Value = '  +2.75%  '
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.92'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.73'
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  +1.98%  '
$ws.Range('D33').Value = '1.427.31'
$ws.Range('E33').Value = '  -2.58%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.10'
$ws.Range('E34').Value = '  +1.14%  '
$ws.Range('E35').Value = '  +1.63%  '
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.891'
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('E38').Value = '  +0.70%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.559'
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.923'
$ws.Range('E40').Value = '  -2.07%  '
$ws.Range('E41').Value = '  +2.24%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '68.67'
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.45'
$ws.Range('E45').Value = '  +3.13%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.82'
$ws.Range('E46').Value = '  +3.29%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '1.784.53'
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '89.29'
$ws.Range('E49').Value = '  +2.09%  '
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('E51').Value = '  +0.57%  '
